$wb = $excel.ActiveWorkbook

# --- Sheet "Replacement": remove the I:L "header" text (row 1) while
#     clearing out the data values underneath it (rows 2-4), keeping the
#     existing number formatting on those data cells. ---
$wsRepl = $wb.Worksheets.Item("Replacement")
$wsRepl.Range("I1:L1").Clear()
$wsRepl.Range("I2:L4").ClearContents()

# Make "Replacement" the active/selected sheet, scroll right so column G
# becomes the left-most visible column, and select I1:L4 with I1 active.
$wsRepl.Activate()
$wsRepl.Range("I1:L4").Select()
$excel.ActiveWindow.ScrollColumn = 7

# --- Sheet "Inspections" is no longer the selected/active sheet (handled
#     above by activating "Replacement" instead). ---
